$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 875
$ws1.Range("F3").Value = 1817
$ws1.Range("F12").Value = 674
$ws1.Range("F14").Value = 4118
$ws1.Range("F17").Value = 3222
$ws1.Range("F21").Value = 205
$ws1.Range("F22").Value = 2175
$ws1.Range("F25").Value = 2003
$ws1.Range("F26").Value = 429
$ws1.Range("F29").Value = 8992
$ws1.Range("F30").Value = 5819
$ws1.Range("F31").Value = 366
$ws1.Range("F32").Value = 188
$ws1.Range("F40").Value = 71
$ws1.Range("F43").Value = 4700
$ws1.Range("F45").Value = 911

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 8476
$ws3.Range("F3").Value = 382
$ws3.Range("F4").Value = 1393

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 875
$ws4.Range("F3").Value = 382
$ws4.Range("F4").Value = 1393
$ws4.Range("F11").Value = 4118
$ws4.Range("F14").Value = 3222
$ws4.Range("F18").Value = 2175
$ws4.Range("F25").Value = 429
$ws4.Range("F28").Value = 8992
$ws4.Range("F31").Value = 366
$ws4.Range("F32").Value = 188
$ws4.Range("F37").Value = 71
$ws4.Range("F41").Value = 4700
$ws4.Range("F43").Value = 911
